$wb = $excel.ActiveWorkbook

# Start from the existing "Spain" sheet - it already has the exact layout /
# formatting / merged cells / styles that the new "Turkey" template needs -
# and duplicate it immediately after itself (this is how the new sheet9
# content, with the same cell styles/merges as Spain, was produced).
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy([System.Reflection.Missing]::Value, $spain) | Out-Null

# The duplicate lands right after Spain and becomes the active sheet/tab.
$turkey = $wb.Worksheets.Item($spain.Index + 1)
$turkey.Name = "Turkey"

# Market name / part number specific to the new Turkey template (these
# become new shared-string entries: "Turkey Market" / "NGC-3191/T3291").
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3291"

# Column D was resized on the new sheet.
$turkey.Columns.Item(4).ColumnWidth = 17

# Selection on the freshly added sheet.
$turkey.Range("G15").Select() | Out-Null

# Spain keeps a plain A1:D12 selection and is no longer the active tab now
# that Turkey has been added (Turkey is last in the tab order and active).
$spain.Range("A1:D12").Select() | Out-Null

$turkey.Activate() | Out-Null
